# Insert a new "mu" column before the existing "isSelected" column (column I).
# This shifts isSelected, bandwidth, transRate, uploadTime, totalTime one
# column to the right (I->J, J->K, K->L, L->M, M->N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I; this shifts everything from I onward right by one.
$ws.Range("I1").EntireColumn.Insert()

# Header for the new column. EntireColumn.Insert() already carries over the
# formatting (bold font, border, centered alignment) from the old column I,
# matching the other header cells in row 1.
$ws.Range("I1").Value = "mu"

# New "mu" values per row (row 2 through 21).
$muValues = @{
    2  = 246943.4382448861
    3  = 280644.209814786
    4  = 222698.3339546729
    5  = 230422.9937791568
    6  = 235818.7246096765
    7  = 285744.5977404236
    8  = 248196.9201041696
    9  = 267164.2500133727
    10 = 287421.8807236081
    11 = 219275.9008517889
    12 = 268456.969184203
    13 = 221636.7964397329
    14 = 152286.8420783606
    15 = 278176.3779967169
    16 = 247711.3304067207
    17 = 246548.3534563272
    18 = 270874.8185469407
    19 = 262031.8998925697
    20 = 246805.3965717752
    21 = 227410.1733562206
}

foreach ($row in $muValues.Keys) {
    $ws.Cells.Item($row, 9).Value = $muValues[$row]
}
